# Auto - Update data with bot!
# Applies updates to the "title" (D) and "link" (E) columns for several
# rows in the blog list worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "모바일 서비스를 위한 다양한 지표"
$ws.Range("E3").Value = "https://lumiamitie.github.io/data/various-metrics-mobile-service/"

$ws.Range("D6").Value = "[Markdown] Jupyter notebook tab like r markdown(widgets)"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Markdown-Jupyter-notebook-tab-like-r-markdownwidgets"

$ws.Range("D9").Value = "Conditional offer 이신 분들한테 연락드립니다 – 2"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/notice-to-conditional-offers-2/#utm_source=rss&utm_medium=rss&utm_campaign=notice-to-conditional-offers-2"

$ws.Range("D12").Value = "사이킷런 1.0 RC1이 릴리즈되었습니다."
$ws.Range("E12").Value = "https://tensorflow.blog/2021/09/13/%ec%82%ac%ec%9d%b4%ed%82%b7%eb%9f%b0-1-0-rc1%ec%9d%b4-%eb%a6%b4%eb%a6%ac%ec%a6%88%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

$ws.Range("D16").Value = "[백준11726, 다이나믹 프로그래밍] 2 x n 타일링 - Python"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/182"

$ws.Range("D37").Value = "[Paper Review] The Origins and Prevalence of Texture Bias in Convolutional Neural Networks"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1820&mod=document&pageid=1"

$ws.Range("D41").Value = "로그 데이터의 수집과 시각화 – Part 2"
$ws.Range("E41").Value = "http://cloudinsight.net/data/log-part-2/"
